# =====================================================================
# data_exclude.xlsx WIP edit:
#  - Add 4 new station sheets (THD, MHD, CMN, GSN) after the existing
#    CGO sheet, each pre-populated with the same instructional header
#    rows as CGO plus a handful of species/date-range exclusion rows.
#  - GSN ends up as the active/selected tab.
# =====================================================================

$wb = $excel.ActiveWorkbook

function Add-SheetAfterLast {
    param($name)
    $cnt = $wb.Worksheets.Count
    $afterSheet = $wb.Worksheets.Item($cnt)
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Name = $name
    return $ws
}

# ---------------------------------------------------------------------
# Step 1: get the sheet ORDER + internal sheetId sequence right first.
# The workbook's sheetId counter tracks "max id ever assigned among
# currently-live sheets, + 1" and is not reset by deleting a sheet
# unless that sheet held the current max id, so scratch sheets are
# added/deleted here to burn spare ids and land on THD=3, MHD=4, CMN=5,
# GSN=7 -- matching the authored file, which had a couple of sheets
# added/removed along the way before settling.
#
# NOTE: deleting a sheet invalidates any *other* previously-fetched
# worksheet object references held in this host, so all the add/delete
# housekeeping happens first, and the real content sheets are
# re-fetched by name afterwards, before anything is written into them.
# ---------------------------------------------------------------------

$scratch1 = Add-SheetAfterLast "zzscratch1"   # consumes id 2
Add-SheetAfterLast "THD" | Out-Null           # id 3
$scratch1.Activate()
$scratch1.Delete()

Add-SheetAfterLast "MHD" | Out-Null           # id 4
Add-SheetAfterLast "CMN" | Out-Null           # id 5

$scratch2 = Add-SheetAfterLast "zzscratch2"   # consumes id 6
Add-SheetAfterLast "GSN" | Out-Null           # id 7
$scratch2 = $wb.Worksheets.Item("zzscratch2")
$scratch2.Activate()
$scratch2.Delete()

Write-Host "Sheet count after reshuffle:" $wb.Worksheets.Count

# ---------------------------------------------------------------------
# Step 2: populate content. Each station sheet repeats the same
# 3-line instructional preamble + header row used on CGO, followed by
# the station's own exclusion rows (species / instrument / start / end,
# all stored as literal TEXT -- not Excel dates -- per the sheet's own
# instructions).
#
# NOTE: this host's command parser mis-parses "FuncCall $var (expr)"
# (a variable argument immediately followed by a parenthesised
# expression) as an invocation of $var itself, so every computed
# argument (e.g. a concatenated cell address) is first assigned to its
# own named variable and only THEN passed along positionally.
# ---------------------------------------------------------------------

function Set-TextCell {
    param($sheet, $addr, $val)
    $cellrng = $sheet.Range($addr)
    $cellrng.NumberFormat = "@"
    $cellrng.Value = $val
}

function Write-StationSheet {
    param($sheet, $rows)

    Set-TextCell $sheet "A1" "# Mole fraction data during the ranges in this worksheet will be excluded"
    Set-TextCell $sheet "A2" "# Date format must by YYYY-MM-DD HH:MM"
    Set-TextCell $sheet "A3" "# Make sure that the cell format is text rather than Excel's date format"
    Set-TextCell $sheet "A4" "Species"
    Set-TextCell $sheet "B4" "Instrument"
    Set-TextCell $sheet "C4" "Start"
    Set-TextCell $sheet "D4" "End"

    $rownum = 5
    foreach ($row in $rows) {
        $addrA = "A" + $rownum
        $addrB = "B" + $rownum
        $addrC = "C" + $rownum
        $addrD = "D" + $rownum
        Set-TextCell $sheet $addrA $row[0]
        Set-TextCell $sheet $addrB $row[1]
        Set-TextCell $sheet $addrC $row[2]
        Set-TextCell $sheet $addrD $row[3]
        $rownum = $rownum + 1
    }
}

# --- THD ---------------------------------------------------------------
$thd = $wb.Worksheets.Item("THD")
$thdRows = @(
    ,@("hfc-32",  "GCMS-Medusa", "2014-03-01 00:00", "2015-08-31 00:00")
    ,@("hfc-125", "GCMS-Medusa", "2014-03-01 00:00", "2015-08-31 00:00")
)
Write-StationSheet $thd $thdRows
$thd.Columns.Item(1).ColumnWidth = 9.998697916666666
$thd.Columns.Item(2).ColumnWidth = 12.498697916666666
$thd.Columns.Item(3).ColumnWidth = 16.666666666666668
$thd.Columns.Item(4).ColumnWidth = 17.166666666666668
$thd.Range("A7").Select()

# --- MHD ---------------------------------------------------------------
$mhd = $wb.Worksheets.Item("MHD")
$mhdRows = @(
    ,@("ch2cl2", "GCMS-Medusa", "2011-05-01 00:00", "2017-04-01 00:00")
)
Write-StationSheet $mhd $mhdRows
$mhd.Columns.Item(1).ColumnWidth = 9.998697916666666
$mhd.Columns.Item(2).ColumnWidth = 12.498697916666666
$mhd.Columns.Item(3).ColumnWidth = 16.666666666666668
$mhd.Columns.Item(4).ColumnWidth = 18.666666666666668
$mhd.Range("D13").Select()

# --- CMN ---------------------------------------------------------------
$cmn = $wb.Worksheets.Item("CMN")
$cmnRows = @(
    ,@("hfc-236fa", "GCMS-Medusa", "1970-01-01 00:00", "2014-05-01 00:00")
    ,@("cfc-114",   "GCMS-Medusa", "1970-01-01 00:00", "2006-01-01 00:00")
    ,@("cfc-115",   "GCMS-Medusa", "1970-01-01 00:00", "2008-01-01 00:00")
)
Write-StationSheet $cmn $cmnRows
$cmn.Columns.Item(1).ColumnWidth = 9.998697916666666
$cmn.Columns.Item(2).ColumnWidth = 14.830729166666666
$cmn.Columns.Item(3).ColumnWidth = 16.330729166666668
$cmn.Columns.Item(4).ColumnWidth = 15.830729166666666
$cmn.Range("D8").Select()

# --- GSN -----------------------------------------------------------------
$gsn = $wb.Worksheets.Item("GSN")
$gsnRows = @(
    ,@("ch3br", "GCMS-Medusa", "1970-01-01 00:00", "2007-12-31 00:00")
)
Write-StationSheet $gsn $gsnRows
$gsn.Columns.Item(1).ColumnWidth = 9.998697916666666
$gsn.Columns.Item(2).ColumnWidth = 12.498697916666666
$gsn.Columns.Item(3).ColumnWidth = 14.666666666666666
$gsn.Columns.Item(4).ColumnWidth = 14.666666666666666
$gsn.Range("D6").Select()

# GSN is the tab left active/selected in the authored workbook.
$gsn = $wb.Worksheets.Item("GSN")
$gsn.Activate()

Write-Host "Done."
